# database/industries/palayesh/shebandar/income/yearly/rial.xlsx
# "update database and change read_price algorithm"
#
# Shift the 5-year reporting window forward by one year (drop 1396/12,
# add 1401/12) together with its matching "publish date" column, and
# refresh the financial figures coming out of the (buggy) updated
# read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "12 ماهه منتهی به ..." period headers (shift one year forward) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish-date headers (shift one year forward) ---
$ws.Range("D9").Value = "1399-04-22 (11)"
$ws.Range("E9").Value = "1400-04-26 (13)"
$ws.Range("F9").Value = "1401-04-27 (10)"
$ws.Range("G9").Value = "1402-02-25 (8)"
$ws.Range("H9").Value = "1402-02-25"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 9

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -4
$ws.Range("E12").Value = -5
$ws.Range("F12").Value = -3
$ws.Range("G12").Value = -7
$ws.Range("H12").Value = -8

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (D15 already "-"; fill E:H with "-") ---
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (now "-" across) ---
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1

# --- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ---
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS on latest capital) ---
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
